# feat: add student import example
# Populate the first two data rows of the "Student" table with a sample
# regular-programme student and a sample international-programme student,
# illustrating how the import template should be filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student")

# Row 2 - Thai-named regular student example
$ws.Range("A2").Value = 63070501000
$ws.Range("B2").Value = "ชื่อ"
$ws.Range("C2").Value = "นามกสุล"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "regular"
$ws.Range("F2").Value = "computer engineer"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "ไม่ระบุ"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 2566
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = "-"

# Row 3 - English-named international student example
$ws.Range("A3").Value = 63070503400
$ws.Range("B3").Value = "name"
$ws.Range("C3").Value = "lastname"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "international"
$ws.Range("F3").Value = "computer engineer"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "ไม่ระบุ"
$ws.Range("L3").Value = "-"
$ws.Range("M3").Value = 2566
$ws.Range("N3").Value = "-"
$ws.Range("O3").Value = "-"
